# Auto-generated edit script applying the Coeurl_Profits.xlsx value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 742845.4
$ws.Range("J17").Value = 742845.4
$ws.Range("L17").Value = 2228536.2
$ws.Range("N17").Value = -2228872.2
$ws.Range("H43").Value = 1150.875
$ws.Range("I43").Value = 1153.3334
$ws.Range("J43").Value = 1149.4
$ws.Range("K43").Value = 1153.3334
$ws.Range("L43").Value = 1149.4
$ws.Range("M43").Value = -1084.3334
$ws.Range("N43").Value = -1287.4
$ws.Range("H132").Value = 1832.4412
$ws.Range("I132").Value = 1377.12
$ws.Range("J132").Value = 3097.2222
$ws.Range("K132").Value = 4131.36
$ws.Range("L132").Value = 9291.6666
$ws.Range("M132").Value = -1601.36
$ws.Range("N132").Value = -14351.6666
$ws.Range("H133").Value = 77900
$ws.Range("J133").Value = 77900
$ws.Range("L133").Value = 77900
$ws.Range("N133").Value = -88020
$ws.Range("H137").Value = 1727.5
$ws.Range("I137").Value = 1555
$ws.Range("J137").Value = 1900
$ws.Range("K137").Value = 4665
$ws.Range("L137").Value = 5700
$ws.Range("N137").Value = -10800
$ws.Range("M137").Value = -2115

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4686.41
$ws.Range("I32").Value = 4406.054
$ws.Range("K32").Value = 4406.054
$ws.Range("M32").Value = -4119.054
$ws.Range("H53").Value = 17398.143
$ws.Range("J53").Value = 20000
$ws.Range("L53").Value = 20000
$ws.Range("N53").Value = -21364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 87416.5
$ws.Range("J132").Value = 87416.5
$ws.Range("L132").Value = 87416.5
$ws.Range("N132").Value = -97536.5
$ws.Range("H140").Value = 96849.5
$ws.Range("J140").Value = 96849.5
$ws.Range("L140").Value = 96849.5
$ws.Range("N140").Value = -107209.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21592.691
$ws.Range("I31").Value = 31139.795
$ws.Range("J31").Value = 3559.2778
$ws.Range("K31").Value = 31139.795
$ws.Range("L31").Value = 3559.2778
$ws.Range("M31").Value = -30844.795
$ws.Range("N31").Value = -4149.2778
$ws.Range("H34").Value = 21592.691
$ws.Range("I34").Value = 31139.795
$ws.Range("J34").Value = 3559.2778
$ws.Range("K34").Value = 31139.795
$ws.Range("L34").Value = 3559.2778
$ws.Range("M34").Value = -30937.795
$ws.Range("N34").Value = -3963.2778
$ws.Range("H103").Value = 9000
$ws.Range("I103").Value = 9000
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 9000
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -7828
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 709.36365
$ws.Range("J55").Value = 2301.6667
$ws.Range("L55").Value = 6905.000100000001
$ws.Range("N55").Value = -7259.000100000001
$ws.Range("H114").Value = 751.9091
$ws.Range("I114").Value = 716.6
$ws.Range("J114").Value = 1105
$ws.Range("K114").Value = 2149.8
$ws.Range("L114").Value = 3315
$ws.Range("M114").Value = 1104.2
$ws.Range("N114").Value = -9823
$ws.Range("H140").Value = 2828.2727
$ws.Range("I140").Value = 2861.2
$ws.Range("K140").Value = 8583.599999999999
$ws.Range("M140").Value = -3403.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 683.7619
$ws.Range("I107").Value = 825.6923
$ws.Range("J107").Value = 453.125
$ws.Range("K107").Value = 825.6923
$ws.Range("L107").Value = 453.125
$ws.Range("M107").Value = 1094.3077
$ws.Range("N107").Value = -4293.125
$ws.Range("H113").Value = 1979.5
$ws.Range("I113").Value = 1977.2222
$ws.Range("K113").Value = 1977.2222
$ws.Range("M113").Value = 192.7778000000001
$ws.Range("H136").Value = 39630
$ws.Range("J136").Value = 39630
$ws.Range("L136").Value = 118890
$ws.Range("N136").Value = -123990
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7668.148
$ws.Range("I7").Value = 8302.75
$ws.Range("J7").Value = 5855
$ws.Range("K7").Value = 8302.75
$ws.Range("L7").Value = 5855
$ws.Range("M7").Value = -8190.75
$ws.Range("N7").Value = -6079
$ws.Range("H12").Value = 732.125
$ws.Range("I12").Value = 601
$ws.Range("K12").Value = 601
$ws.Range("M12").Value = -431
$ws.Range("H40").Value = 6315.6665
$ws.Range("I40").Value = 5778.8
$ws.Range("K40").Value = 5778.8
$ws.Range("M40").Value = -5642.8
$ws.Range("H46").Value = 1533
$ws.Range("I46").Value = 1299.5
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1299.5
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -1111.5
$ws.Range("N46").Value = -2376
$ws.Range("H108").Value = 39900
$ws.Range("J108").Value = 39900
$ws.Range("L108").Value = 39900
$ws.Range("N108").Value = -47580
$ws.Range("H126").Value = 7668.148
$ws.Range("I126").Value = 8302.75
$ws.Range("J126").Value = 5855
$ws.Range("K126").Value = 24908.25
$ws.Range("L126").Value = 17565
$ws.Range("M126").Value = -22438.25
$ws.Range("N126").Value = -22505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 6479.25
$ws.Range("I9").Value = 8529
$ws.Range("J9").Value = 330
$ws.Range("K9").Value = 8529
$ws.Range("L9").Value = 330
$ws.Range("M9").Value = -8389
$ws.Range("N9").Value = -610
$ws.Range("H113").Value = 1028.3928
$ws.Range("I113").Value = 718.3125
$ws.Range("J113").Value = 1441.8334
$ws.Range("K113").Value = 2154.9375
$ws.Range("L113").Value = 4325.5002
$ws.Range("M113").Value = 15.0625
$ws.Range("N113").Value = -8665.5002
$ws.Range("H132").Value = 1425.762
$ws.Range("I132").Value = 1372.05
$ws.Range("K132").Value = 4116.15
$ws.Range("M132").Value = -1586.15
$ws.Range("H136").Value = 1782.6735
$ws.Range("I136").Value = 1422
$ws.Range("J136").Value = 2461.5881
$ws.Range("K136").Value = 4266
$ws.Range("L136").Value = 7384.7643
$ws.Range("M136").Value = -1716
$ws.Range("N136").Value = -12484.7643
$ws.Range("H140").Value = 74697.25
$ws.Range("J140").Value = 74697.25
$ws.Range("L140").Value = 74697.25
$ws.Range("N140").Value = -85057.25
